$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value = "A's Astounding Team"
$ws.Range("B10").Value = "Magic M"
$ws.Range("B12").Value = "K's Deluxe Team"

$ws.Range("B13").Select()
